$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row (row 1) with new sequence values in P1, Q1,
# copying the style/format used by the rest of row 1 (bold header style).
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

# Update columns I, K, M, O for data rows 2-25, and populate new
# columns P and Q with value 2.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P -> 2 (new)
    $ws.Cells.Item($r, 17).Value = 2   # Q -> 2 (new)
}
